$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 68, shifting existing rows 68-99 down to 69-100.
$ws.Rows.Item(68).Insert()

# Populate the new row 68 with the new record.
$ws.Cells.Item(68, 1).Value = 10
$ws.Cells.Item(68, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(68, 3).Value = "La Araucanía"
$ws.Cells.Item(68, 4).Value = 44875
$ws.Cells.Item(68, 5).Value = 9
$ws.Cells.Item(68, 6).Value = 100112022
$ws.Cells.Item(68, 7).Value = "Arveja Verde"
$ws.Cells.Item(68, 8).Value = "Sin especificar"
$ws.Cells.Item(68, 9).Value = "Primera"
$ws.Cells.Item(68, 10).Value = 35
$ws.Cells.Item(68, 11).Value = 25000
$ws.Cells.Item(68, 12).Value = 25000
$ws.Cells.Item(68, 13).Value = 25000
$ws.Cells.Item(68, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(68, 15).Value = "Región Metropolitana"
$ws.Cells.Item(68, 16).Value = 1000
$ws.Cells.Item(68, 17).Value = 25
$ws.Cells.Item(68, 18).Value = "Hortaliza"

# Match the date-format style used by the other rows' Fecha (column D) cells.
$ws.Cells.Item(68, 4).NumberFormat = $ws.Cells.Item(69, 4).NumberFormat
